$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (nombre_aides, montant_total)
$updates = @{
    13  = @(187863, 1168289738)
    88  = @(71283,  110329263)
    91  = @(18884,  75370883)
    93  = @(16953,  50862471)
    98  = @(6301,   19502085)
    100 = @(9348,   23899308)
    121 = @(1306416, 2275552863)
    122 = @(382,    1260995)
    129 = @(633827, 3435576311)
    130 = @(4250,   141666559)
    132 = @(586042, 3473327690)
    136 = @(26706,  144413770)
    178 = @(515892, 891219559)
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("C$row").Value = $values[0]
    $ws.Range("E$row").Value = $values[1]
}
